$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K becomes F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting (date / integer styles) from column F into the
# newly inserted D:E columns so the new quarter columns render the same
# way as the existing ones.
$ws.Range("F1:F102").Copy()
$ws.Range("D1:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 3691000
$ws.Range("E8").Value = 3729000
$ws.Range("D9").Value = 3204000
$ws.Range("E9").Value = 3172000
$ws.Range("F9").Value = 3120000
$ws.Range("H9").Value = 2748000
$ws.Range("I9").Value = 2828000
$ws.Range("D10").Value = 487000
$ws.Range("E10").Value = 557000
$ws.Range("F10").Value = 489000
$ws.Range("H10").Value = 385000
$ws.Range("I10").Value = 420000
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 3000
$ws.Range("H14").Value = 42000
$ws.Range("D15").Value = 137000
$ws.Range("E15").Value = 126000
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("D17").Value = 3399000
$ws.Range("E17").Value = 3359000
$ws.Range("H17").Value = 2993000
$ws.Range("I17").Value = 3019000
$ws.Range("D18").Value = 292000
$ws.Range("E18").Value = 370000
$ws.Range("H18").Value = 140000
$ws.Range("I18").Value = 229000
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("D20").Value = -5000
$ws.Range("E20").Value = -15000
$ws.Range("H20").Value = -17000
$ws.Range("I20").Value = -22000
$ws.Range("D21").Value = 424000
$ws.Range("E21").Value = 481000
$ws.Range("D22").Value = 34000
$ws.Range("E22").Value = 41000
$ws.Range("D23").Value = 253000
$ws.Range("E23").Value = 314000
$ws.Range("D24").Value = -339000
$ws.Range("E24").Value = 23000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 592000
$ws.Range("E26").Value = 291000
$ws.Range("D27").Value = 592000
$ws.Range("E27").Value = 291000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 5000
$ws.Range("E32").Value = 15000
$ws.Range("H32").Value = 17000
$ws.Range("I32").Value = 22000
$ws.Range("D33").Value = 592000
$ws.Range("E33").Value = 291000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 592000
$ws.Range("E35").Value = 291000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("D41").Value = 1000000
$ws.Range("E41").Value = 1344000
$ws.Range("D42").Value = 12000
$ws.Range("E42").Value = 12000
$ws.Range("F42").Value = 11000
$ws.Range("G42").Value = "NA"
$ws.Range("H42").Value = "NA"
$ws.Range("I42").Value = "NA"
$ws.Range("J42").Value = "NA"
$ws.Range("D43").Value = 1645000
$ws.Range("E43").Value = 1660000
$ws.Range("F43").Value = 1644000
$ws.Range("D44").Value = 2092000
$ws.Range("E44").Value = 1950000
$ws.Range("D45").Value = 81000
$ws.Range("E45").Value = 102000
$ws.Range("D46").Value = 4830000
$ws.Range("E46").Value = 5068000
$ws.Range("D47").Value = 513000
$ws.Range("E47").Value = 508000
$ws.Range("D48").Value = 4865000
$ws.Range("E48").Value = 4643000
$ws.Range("D49").Value = 158000
$ws.Range("E49").Value = 160000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 616000
$ws.Range("E52").Value = 190000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 10982000
$ws.Range("E54").Value = 10569000
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("D57").Value = 2517000
$ws.Range("E57").Value = 2510000
$ws.Range("D58").Value = 65000
$ws.Range("E58").Value = 4000
$ws.Range("D59").Value = 615000
$ws.Range("E59").Value = 622000
$ws.Range("D60").Value = 3197000
$ws.Range("E60").Value = 3136000
$ws.Range("D61").Value = 2316000
$ws.Range("E61").Value = 2498000
$ws.Range("D62").Value = 1266000
$ws.Range("E62").Value = 993000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6780000
$ws.Range("E66").Value = 6628000
$ws.Range("D67").Value = ""
$ws.Range("E67").Value = ""
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1212000
$ws.Range("E72").Value = 629000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 4202000
$ws.Range("E76").Value = 3941000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 592000
$ws.Range("E81").Value = 291000
$ws.Range("D82").Value = ""
$ws.Range("E82").Value = ""
$ws.Range("D83").Value = 137000
$ws.Range("E83").Value = 126000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 216000
$ws.Range("E89").Value = 429000
$ws.Range("H89").Value = 280000
$ws.Range("I89").Value = 303000
$ws.Range("D90").Value = ""
$ws.Range("E90").Value = ""
$ws.Range("D91").Value = -355000
$ws.Range("E91").Value = -265000
$ws.Range("I91").Value = -171000
$ws.Range("J91").Value = -73000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -326000
$ws.Range("E94").Value = -256000
$ws.Range("H94").Value = -197000
$ws.Range("I94").Value = -69000
$ws.Range("D95").Value = ""
$ws.Range("E95").Value = ""
$ws.Range("D96").Value = -9000
$ws.Range("E96").Value = -9000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -231000
$ws.Range("E100").Value = -57000
$ws.Range("D101").Value = -4000
$ws.Range("E101").Value = -3000
$ws.Range("D102").Value = -345000
$ws.Range("E102").Value = 113000
$ws.Range("H102").Value = -139000
$ws.Range("I102").Value = 171000

Write-Output "done"
